$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rows 29, 30 and 35 had their match data cyclically rotated:
#    new row29 <- old row30, new row30 <- old row35, new row35 <- old row29
#    (the columns that stay identical across the three rows - G, K, O, S -
#     are left untouched)
# ---------------------------------------------------------------------------

$ws.Range("F29").Value = "Charlton"
$ws.Range("H29").Value = "Bristol Rovers"
$ws.Range("I29").Value = 2
$ws.Range("J29").Value = 2.07
$ws.Range("L29").Value = 2.28
$ws.Range("M29").Value = "15/08/2023 20:43"
$ws.Range("N29").Value = 3.57
$ws.Range("P29").Value = 3.55
$ws.Range("Q29").Value = "15/08/2023 20:41"
$ws.Range("R29").Value = 3.56
$ws.Range("T29").Value = 3.18
$ws.Range("U29").Value = "15/08/2023 20:43"
$ws.Range("V29").Value = "https://www.betexplorer.com/football/england/league-one/charlton-bristol-rovers/SrvESH0P/"

$ws.Range("F30").Value = "Carlisle"
$ws.Range("H30").Value = "Wigan"
$ws.Range("I30").Value = 1
$ws.Range("J30").Value = 2.68
$ws.Range("L30").Value = 2.71
$ws.Range("M30").Value = "15/08/2023 20:02"
$ws.Range("N30").Value = 3.26
$ws.Range("P30").Value = 3.24
$ws.Range("Q30").Value = "15/08/2023 18:42"
$ws.Range("R30").Value = 2.76
$ws.Range("T30").Value = 2.81
$ws.Range("U30").Value = "15/08/2023 20:02"
$ws.Range("V30").Value = "https://www.betexplorer.com/football/england/league-one/carlisle-wigan/EujW3YNO/"

$ws.Range("F35").Value = "Derby"
$ws.Range("H35").Value = "Oxford Utd"
$ws.Range("I35").Value = 2
$ws.Range("J35").Value = 1.85
$ws.Range("L35").Value = 1.96
$ws.Range("M35").Value = "15/08/2023 19:25"
$ws.Range("N35").Value = 3.59
$ws.Range("P35").Value = 3.43
$ws.Range("Q35").Value = "15/08/2023 19:25"
$ws.Range("R35").Value = 4.41
$ws.Range("T35").Value = 4.27
$ws.Range("U35").Value = "15/08/2023 19:25"
$ws.Range("V35").Value = "https://www.betexplorer.com/football/england/league-one/derby-oxford-utd/44VuOy9t/"

# ---------------------------------------------------------------------------
# 2) Append three new match rows (88, 89, 90), pulling the A/E cell
#    formatting from the last existing data row (87) so the new rows match
#    the same styling (bordered index column, date-formatted match-date
#    column).
# ---------------------------------------------------------------------------

$ws.Range("A87:V87").Copy()
$ws.Range("A88:V90").PasteSpecial(-4122)  # xlPasteFormats

$newRows = @(
    @{ Row = 88; A = 87; F = "Barnsley";     G = 2; H = "Portsmouth"; I = 3;
       J = 2.15; K = "16/09/2023 17:13"; L = 2.43; M = "19/09/2023 20:44";
       N = 3.64; O = "16/09/2023 17:13"; P = 3.52; Q = "19/09/2023 20:44";
       R = 3.12; S = "16/09/2023 17:13"; T = 2.96; U = "19/09/2023 20:44";
       V = "https://www.betexplorer.com/football/england/league-one/barnsley-portsmouth/rDiXc0dR/" },
    @{ Row = 89; A = 88; F = "Peterborough"; G = 3; H = "Cheltenham";  I = 0;
       J = 1.5;  K = "16/09/2023 17:13"; L = 1.47; M = "19/09/2023 20:41";
       N = 4.39; O = "16/09/2023 17:13"; P = 4.61; Q = "19/09/2023 20:41";
       R = 5.78; S = "16/09/2023 17:13"; T = 7;    U = "19/09/2023 20:41";
       V = "https://www.betexplorer.com/football/england/league-one/peterborough-cheltenham/UXVFkx47/" },
    @{ Row = 90; A = 89; F = "Port Vale";    G = 2; H = "Burton";      I = 3;
       J = 1.88; K = "16/09/2023 17:13"; L = 2.04; M = "19/09/2023 20:41";
       N = 3.59; O = "16/09/2023 17:13"; P = 3.42; Q = "19/09/2023 20:41";
       R = 4.28; S = "16/09/2023 17:13"; T = 3.97; U = "19/09/2023 20:41";
       V = "https://www.betexplorer.com/football/england/league-one/port-vale-burton/ryzKldJD/" }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Range("A$row").Value = $r.A
    $ws.Range("B$row").Value = "england"
    $ws.Range("C$row").Value = "league-one"
    $ws.Range("D$row").Value = "2023-2024"
    $ws.Range("E$row").Value = 45188.86458333334
    $ws.Range("F$row").Value = $r.F
    $ws.Range("G$row").Value = $r.G
    $ws.Range("H$row").Value = $r.H
    $ws.Range("I$row").Value = $r.I
    $ws.Range("J$row").Value = $r.J
    $ws.Range("K$row").Value = $r.K
    $ws.Range("L$row").Value = $r.L
    $ws.Range("M$row").Value = $r.M
    $ws.Range("N$row").Value = $r.N
    $ws.Range("O$row").Value = $r.O
    $ws.Range("P$row").Value = $r.P
    $ws.Range("Q$row").Value = $r.Q
    $ws.Range("R$row").Value = $r.R
    $ws.Range("S$row").Value = $r.S
    $ws.Range("T$row").Value = $r.T
    $ws.Range("U$row").Value = $r.U
    $ws.Range("V$row").Value = $r.V
}
